$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.262.20"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "3.815.46"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "706.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "3.813.40"
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("E11").Value = "  +5.39%  "
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "4.457.76"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").Value = "3.801.61"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "71.220.55"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "513.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.724"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000145"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("D26").Value = "3.965.67"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("E30").Value = "  -3.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("E35").Value = "  -5.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").Value = "3.775.87"
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("E42").Value = "  -1.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "170.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.10%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "423.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.295"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.45%  "
